# Auto-update predictions and index for 2025-10-18
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: AS Saint-Étienne - Le Mans FC
$ws.Range("A2").Value = "AS Saint-Étienne  - Le Mans FC: 19:00"
$ws.Range("B2").Value = "AS Saint-Étienne"
$ws.Range("C2").Value = 71
$ws.Range("D2").Value = 100
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = 1.45
$ws.Range("G2").Value = ""

# Row 3: UE Santa Coloma B - Sporting Club Escaldes
$ws.Range("A3").Value = "UE Santa Coloma B - Sporting Club Escaldes : 19:30"
$ws.Range("B3").Value = "Sporting Club Escaldes"
$ws.Range("C3").Value = 70
$ws.Range("D3").Value = 86
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = 1.67
$ws.Range("G3").Value = ""

# Row 4: Västerås SK - GIF Sundsvall (result known: 2:1, pick correct)
$ws.Range("A4").Value = "Västerås SK ✓ - GIF Sundsvall: 2:1"
$ws.Range("B4").Value = "Västerås SK"
$ws.Range("C4").Value = 70
$ws.Range("D4").Value = 85
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = 1.62
$ws.Range("G4").Value = "✓"

# Row 5: Thep Xanh Nam Dinh FC - Becamex Ho Chi Minh City FC (result known: 1:2, pick wrong)
$ws.Range("A5").Value = "Thep Xanh Nam Dinh FC X - Becamex Ho Chi Minh City FC: 1:2"
$ws.Range("B5").Value = "Thep Xanh Nam Dinh FC"
$ws.Range("C5").Value = 67
$ws.Range("D5").Value = 75
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = 1.36
$ws.Range("G5").Value = "X"

# Row 6: Al-Ain FC - FC Baniyas (result known: 4:0, pick correct)
$ws.Range("A6").Value = "Al-Ain FC ✓ - FC Baniyas: 4:0"
$ws.Range("B6").Value = "Al-Ain FC"
$ws.Range("C6").Value = 59
$ws.Range("D6").Value = 100
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = 1.45
$ws.Range("G6").Value = "✓"

# Row 7 (new): Club Nacional - Miramar Misiones
$ws.Range("A7").Value = "Club Nacional  - Miramar Misiones: 22:30"
$ws.Range("B7").Value = "Club Nacional"
$ws.Range("C7").Value = 59
$ws.Range("D7").Value = 86
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = 1.67
$ws.Range("G7").Value = ""

# Row 8 (new): Olympique Marseille - Le Havre AC
$ws.Range("A8").Value = "Olympique Marseille  - Le Havre AC: 20:05"
$ws.Range("B8").Value = "Olympique Marseille"
$ws.Range("C8").Value = 58
$ws.Range("D8").Value = 97
$ws.Range("E8").Value = 94
$ws.Range("F8").Value = 1.38
$ws.Range("G8").Value = ""
